$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---
$ws.Range("C2").Value = "Hartmut"

# Card number must stay as text (it would otherwise be parsed as a number)
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 06.08.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "10.08."
$ws.Range("C6").Value = "11.08."
$ws.Range("D6").Value = "BURGER KING Hofgeismar"
$ws.Range("E6").Value = "40,03-"

# --- Row 7 ---
$ws.Range("B7").Value = "14.08."
$ws.Range("C7").Value = "15.08."
$ws.Range("D7").Value = "KARTENZ./14.08 EDEKA RO"
$ws.Range("E7").Value = "109,05-"

# --- Row 8 ---
$ws.Range("B8").Value = "15.08."
$ws.Range("C8").Value = "16.08."
$ws.Range("D8").Value = "PAYPAL TCZFMB"
$ws.Range("E8").Value = "52,10-"

# --- Row 9 ---
$ws.Range("B9").Value = "19.08."
$ws.Range("C9").Value = "20.08."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 52028385"
$ws.Range("E9").Value = "41,70-"

# --- Row 10 ---
$ws.Range("B10").Value = "22.08."
$ws.Range("C10").Value = "23.08."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-14588052"
$ws.Range("E10").Value = "57,24-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 25.08.2024"
$ws.Range("E12").Value = "300,12-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 31.08.2024"
